$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values
$ws.Range("B3").Value = 0.05987764829982773
$ws.Range("C3").Value = 0.6520103082685891
$ws.Range("D3").Value = 0.7338235967626286
$ws.Range("E3").Value = 0.8566350429223805
$ws.Range("F3").Value = 0.8779561813900352
$ws.Range("G3").Value = 19

# Update row 4 values
$ws.Range("B4").Value = 0.251492217976663
$ws.Range("C4").Value = 0.7482559981055826
$ws.Range("D4").Value = 0.8928839205617156
$ws.Range("E4").Value = 0.9449253518462268
$ws.Range("F4").Value = 0.9372500616622482
$ws.Range("G4").Value = 18
